$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.692.22"
$ws.Range("E2").Value = "  +1.36%  "

$ws.Range("D3").Value = "2.305.73"
$ws.Range("E3").Value = "  +0.83%  "

$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.25%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.73%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.607"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.65%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.00"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.88%  "

$ws.Range("E11").Value = "  +0.62%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.52"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.60%  "

$ws.Range("E13").Value = "  +1.00%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.999"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.46%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.42"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.53%  "

$ws.Range("D16").Value = "2.654.72"
$ws.Range("E16").Value = "  +0.74%  "

$ws.Range("D17").Value = "2.304.09"
$ws.Range("E17").Value = "  +0.62%  "

$ws.Range("D18").Value = "42.514.01"
$ws.Range("E18").Value = "  +1.11%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.63"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.33%  "

$ws.Range("E20").Value = "  +0.74%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +33.21%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "74.09"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.21%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.54"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.40%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "268.44"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.86%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.24"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.47%  "

$ws.Range("E26").Value = "  -0.26%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.92"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.56%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.37%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.62"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.70%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.14"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.77%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.56"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +12.73%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "166.23"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.15%  "

$ws.Range("E33").Value = "  +1.90%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.68"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.74%  "

$ws.Range("E35").Value = "  -2.69%  "

$ws.Range("E36").Value = "  -0.43%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.59"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.12%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0354"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.16%  "

$ws.Range("B39").Value = "NEARProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.70"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.43%  "

$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.77"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.90%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.64"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +13.51%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "98.18"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.80%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "70.03"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.26%  "

$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.226"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.98%  "

$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.43%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "118.28"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.50%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.32"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.86%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "79.98"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.31%  "

$ws.Range("D49").Value = "1.649.91"
$ws.Range("E49").Value = "  +4.91%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.31"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.73%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.86"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.23%  "
